$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $r = $ws.Range($cellAddr)
    $r.Value = "'" + $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "321.73"
Set-TextValue "E2" "-3.15%"
Set-TextValue "D3" "42.91"
Set-TextValue "E3" "-6.36%"
Set-TextValue "D4" "5.207"
Set-TextValue "E4" "-7.46%"
Set-TextValue "D5" "0.08187"
Set-TextValue "E5" "-1.97%"
Set-TextValue "D6" "4.320"
Set-TextValue "E6" "-2.94%"
Set-TextValue "D7" "1.831"
Set-TextValue "E7" "-10.31%"
Set-TextValue "D8" "0.9358"
Set-TextValue "E8" "-3.64%"
Set-TextValue "D9" "0.1113"
Set-TextValue "E9" "-4.57%"
Set-TextValue "D10" "0.1862"
Set-TextValue "E10" "-2.99%"
Set-TextValue "D11" "0.09462"
Set-TextValue "E11" "-5.38%"
Set-TextValue "D12" "0.04611"
Set-TextValue "E12" "-1.98%"
Set-TextValue "D13" "7.416"
Set-TextValue "E13" "-28.42%"
Set-TextValue "D14" "0.1057"
Set-TextValue "E14" "-0.21%"
Set-TextValue "D15" "0.001301"
Set-TextValue "E15" "1.00%"
Set-TextValue "D16" "0.005746"
Set-TextValue "E16" "-4.78%"
Set-TextValue "D17" "3.362"
Set-TextValue "E17" "-0.48%"
Set-TextValue "D18" "2.520"
Set-TextValue "E18" "-1.56%"
Set-TextValue "E19" "-0.79%"
Set-TextValue "E20" "-0.39%"
Set-TextValue "E21" "-1.21%"
Set-TextValue "D22" "0.04159"
Set-TextValue "E22" "-0.62%"
Set-TextValue "D23" "0.001248"
Set-TextValue "E23" "-4.82%"
Set-TextValue "D24" "0.004315"
Set-TextValue "E24" "-6.02%"
Set-TextValue "D25" "0.0001100"
Set-TextValue "E25" "-15.44%"
Set-TextValue "E26" "-20.60%"
Set-TextValue "D38" "0.02718"
Set-TextValue "E38" "-1.39%"
Set-TextValue "D39" "0.05554"
Set-TextValue "E39" "-4.19%"
Set-TextValue "D40" "0.007951"
Set-TextValue "E40" "3.33%"
Set-TextValue "D41" "0.1396"
Set-TextValue "E41" "-2.81%"
Set-TextValue "E42" "-10.28%"
Set-TextValue "D43" "0.002093"
Set-TextValue "E43" "4.02%"
Set-TextValue "D44" "0.007510"
Set-TextValue "E44" "-13.91%"
Set-TextValue "D45" "0.3202"
Set-TextValue "E45" "-5.90%"
Set-TextValue "D46" "0.00006971"
Set-TextValue "E46" "-4.33%"
Set-TextValue "E47" "-0.19%"
Set-TextValue "D48" "0.003466"
Set-TextValue "E48" "-1.10%"
Set-TextValue "E49" "0.60%"
Set-TextValue "E50" "-0.19%"
Set-TextValue "E51" "-0.19%"
